# Regenerate column G ("K") values for rows 2-49 on the active sheet.
# These values come from freshly re-scraped/recalculated source data
# (commit message: "regen save_data to use K instead of Strike#, regen
# std/mean, calc and write s_vals"). The sheet stores raw values (no
# formulas), so we simply overwrite the affected cells with the new
# numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 3
    4  = 0
    5  = 0
    6  = 0
    7  = 1
    8  = 0
    9  = 1
    10 = 0
    11 = 1
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 2
    21 = 1
    22 = 0
    23 = 0
    24 = 1
    25 = 1
    26 = 2
    27 = 2
    28 = 1
    29 = 0
    30 = 0
    31 = 0
    32 = 1
    33 = 0
    34 = 2
    35 = 2
    36 = 2
    37 = 1
    38 = 0
    39 = 0
    40 = 0
    41 = 1
    42 = 0
    43 = 1
    44 = 0
    45 = 1
    46 = 0
    47 = 0
    48 = 1
    49 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
